$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the two new mail-log rows (64 and 65) to the "Logs" sheet.
$logs.Range("A64").Value = "Sollicitatie marketingfunctie"
$logs.Range("B64").Value = "mailmind.test@zohomail.eu"
$logs.Range("C64").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D64").Value = "Overig"
$logs.Range("F64").Value = "2025-06-17 23:02:10"
$logs.Range("G64").Value = "Nee"

$logs.Range("A65").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B65").Value = "mailmind.test@zohomail.eu"
$logs.Range("C65").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D65").Value = "Bestelling"
$logs.Range("F65").Value = "2025-06-17 23:02:11"
$logs.Range("G65").Value = "Nee"

# Extend the conditional-formatting ranges so the new rows are covered too.
$logs.Range("D2:D63").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D65"))
$logs.Range("G2:G63").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G65"))

# Update the Dashboard category counts to reflect the new rows.
$dash.Range("B3").Value = 18
$dash.Range("B5").Value = 7
